# June 15 - Inputs update
# Insert a new Kit test-data row ("deluxe25offp-redes-summerb") right after the
# existing "deluxe25offp-redes-spring" row, and rename that existing row's
# campaign value to "deluxe25offp-redes-summera".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 17 (pushes "End" and the repeated table below it down by one row)
$ws.Rows("17:17").Insert()

# Rename the campaign value that was on row 16 and drop its leftover wrap-text style
$ws.Range("C16").Value2 = "deluxe25offp-redes-summera"
$ws.Range("C16").ClearFormats()

# Fill in the newly inserted row 17 with the new test case
$ws.Range("A17").Value2 = "QA"
$ws.Range("B17").Value2 = "Sub-D"
$ws.Range("C17").Value2 = "deluxe25offp-redes-summerb"
$ws.Range("D17").Value2 = "Kit"
$ws.Range("E17").Value2 = "Chrome"
$ws.Range("A17:E17").ClearFormats()

# Match the saved selection shown in the workbook
$ws.Range("A13").Select()
